# Updates cryptos list figures (price/volume) scraped on Sun Dec 10 03:38:21 UTC 2023.
# Rows 10/11 and 43/44/45 also had their Coin/Link swapped (reordered ranking).
#
# Numeric-looking Price strings (single "." decimal, e.g. "0.673") are written
# with a leading apostrophe so Excel keeps them as TEXT (matching the source
# workbook, where every Price/Volume cell is a text cell) instead of auto-
# converting them to numbers. Values with two dots (e.g. "43.817.55") or the
# Volume(1h) percentages (e.g. "  -1.03%  ") are never parsed as numbers by
# Excel, so they need no such protection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '43.817.55'
$ws.Range('E2').Value = '  -1.03%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.351.08'
$ws.Range('E3').Value = '  -0.40%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.04%  '

# Row 5: XRP
$ws.Range('D5').Value = '''0.673'
$ws.Range('E5').Value = '  -1.35%  '

# Row 6: BNB
$ws.Range('D6').Value = '''240.45'

# Row 7: Solana
$ws.Range('D7').Value = '''73.17'
$ws.Range('E7').Value = '  -1.80%  '

# Row 8: USDC
$ws.Range('E8').Value = '  -0.05%  '

# Row 9: Cardano
$ws.Range('D9').Value = '''0.598'
$ws.Range('E9').Value = '  +3.28%  '

# Row 10: OKB (was Dogecoin)
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '''62.05'
$ws.Range('E10').Value = '  +7.66%  '

# Row 11: Dogecoin (was OKB)
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '''0.100'
$ws.Range('E11').Value = '  -2.29%  '

# Row 12: Avalanche
$ws.Range('D12').Value = '''32.76'
$ws.Range('E12').Value = '  +1.59%  '

# Row 13: Polkadot
$ws.Range('D13').Value = '''7.32'
$ws.Range('E13').Value = '  -2.49%  '

# Row 14: TRON
$ws.Range('E14').Value = '  +0.13%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '2.699.32'
$ws.Range('E15').Value = '  -0.51%  '

# Row 16: Chainlink
$ws.Range('D16').Value = '''16.40'
$ws.Range('E16').Value = '  -3.14%  '

# Row 17: Polygon
$ws.Range('D17').Value = '''0.905'
$ws.Range('E17').Value = '  -1.31%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '2.344.99'
$ws.Range('E18').Value = '  -0.54%  '

# Row 19: WrappedBTC
$ws.Range('D19').Value = '43.729.30'
$ws.Range('E19').Value = '  -1.62%  '

# Row 20: ShibaInu
$ws.Range('D20').Value = '''0.0000102'
$ws.Range('E20').Value = '  -1.54%  '

# Row 21: Uniswap
$ws.Range('D21').Value = '''6.69'
$ws.Range('E21').Value = '  -1.54%  '

# Row 22: Litecoin
$ws.Range('D22').Value = '''77.11'
$ws.Range('E22').Value = '  -1.75%  '

# Row 23: BitcoinCash
$ws.Range('D23').Value = '''256.07'
$ws.Range('E23').Value = '  -0.20%  '

# Row 24: ImmutableX
$ws.Range('E24').Value = '  +20.08%  '

# Row 26: WEMIXToken
$ws.Range('E26').Value = '  -0.30%  '

# Row 27: PancakeSwap
$ws.Range('D27').Value = '''2.50'
$ws.Range('E27').Value = '  -3.28%  '

# Row 28: Cosmos
$ws.Range('E28').Value = '  -1.90%  '

# Row 29: Toncoin
$ws.Range('E29').Value = '  +0.75%  '

# Row 30: EthereumClassic
$ws.Range('D30').Value = '''22.62'
$ws.Range('E30').Value = '  +0.01%  '

# Row 31: Monero
$ws.Range('D31').Value = '''177.26'
$ws.Range('E31').Value = '  +1.23%  '

# Row 32: Kaspa
$ws.Range('E32').Value = '  -1.24%  '

# Row 33: Stellar
$ws.Range('E33').Value = '  +1.60%  '

# Row 34: Hedera
$ws.Range('D34').Value = '''0.0755'
$ws.Range('E34').Value = '  -0.33%  '

# Row 35: Filecoin
$ws.Range('E35').Value = '  -5.12%  '

# Row 36: InternetComputer(DFINITY)
$ws.Range('D36').Value = '''5.45'
$ws.Range('E36').Value = '  +1.47%  '

# Row 37: RenderToken
$ws.Range('E37').Value = '  -1.79%  '

# Row 38: LidoDAOToken
$ws.Range('E38').Value = '  -4.11%  '

# Row 39: THORChain
$ws.Range('E39').Value = '  -4.13%  '

# Row 40: VeChain
$ws.Range('E40').Value = '  +1.35%  '

# Row 41: MultiversX
$ws.Range('D41').Value = '''68.80'
$ws.Range('E41').Value = '  +29.60%  '

# Row 42: Cronos
$ws.Range('E42').Value = '  +11.23%  '

# Row 43: FraxShare (was InjectiveProtocol)
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''9.12'
$ws.Range('E43').Value = '  +1.09%  '

# Row 44: FTXToken (was FraxShare)
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').Value = '''4.95'
$ws.Range('E44').Value = '  +10.24%  '

# Row 45: InjectiveProtocol (was FTXToken)
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '''19.09'
$ws.Range('E45').Value = '  -2.04%  '

# Row 46: Algorand
$ws.Range('E46').Value = '  +4.79%  '

# Row 47: NEARProtocol
$ws.Range('D47').Value = '''2.50'
$ws.Range('E47').Value = '  -0.63%  '

# Row 48: BinanceUSD
$ws.Range('E48').Value = '  +0.06%  '

# Row 49: TrustWalletToken
$ws.Range('E49').Value = '  -1.53%  '

# Row 50: ARBITRUM
$ws.Range('E50').Value = '  -1.75%  '

# Row 51: Aave
$ws.Range('D51').Value = '''98.22'
$ws.Range('E51').Value = '  -3.41%  '
